# Update "想去人数" (column F) figures for several events on the
# "展览" and "全部类型" worksheets to reflect newly scraped totals.

$wb = $excel.ActiveWorkbook

# Map of event name -> new "想去人数" value
$updates = @{
    "南宁·小蜜蜂动漫嘉年华2.0" = 284
    "南宁·漫控嘉年华09暨南宁高校动漫联盟十六周年联合漫展" = 97
    "南宁·0713国乙ONLY" = 279
    "广西·首届明日方舟only展 - 花庭圣梦" = 216
    "南宁·AB动漫游戏嘉年华" = 2016
    "南宁·第十九届（2024）良牙动漫夏季盛典（良牙夏典）" = 4785
    "南宁·蔚蓝档案only" = 338
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)

    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count

    for ($r = 2; $r -le $rowCount; $r++) {
        $name = $ws.Cells.Item($r, 3).Value2
        if ($null -ne $name -and $updates.ContainsKey($name)) {
            $ws.Cells.Item($r, 6).Value = $updates[$name]
        }
    }
}
